$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 290, shifting existing row 290 (and everything below)
# down to row 291, growing the used range from A1:R369 to A1:R370.
$ws.Rows("290:290").Insert(-4121)

# Populate the newly inserted row 290 with the new data record.
$ws.Range("A290").Value = 3
$ws.Range("B290").Value = "Femacal de La Calera"
$ws.Range("C290").Value = "Coquimbo"
$ws.Range("D290").Value = 44736
$ws.Range("E290").Value = 5
$ws.Range("F290").Value = 100114013
$ws.Range("G290").Value = "Zanahoria"
$ws.Range("H290").Value = "Sin especificar"
$ws.Range("I290").Value = "Primera"
$ws.Range("J290").Value = 480
$ws.Range("K290").Value = 7000
$ws.Range("L290").Value = 7300
$ws.Range("M290").Value = 7156
$ws.Range("N290").Value = "$/saco 20 kilos"
$ws.Range("O290").Value = "Provincia de Quillota"
$ws.Range("P290").Value = 358
$ws.Range("Q290").Value = 20
$ws.Range("R290").Value = "Hortaliza"
